# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" (fund-level holdings) right before the
#    "总计" (totals) sheet.
# 2. Insert a new top data row into "总计" summarizing the 2022-Q1 quarter,
#    pushing the existing rows down and renumbering the index column.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet with per-fund holdings, inserted before "总计"
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# Worksheets.Add($total) inserted the new sheet *before* "总计", which
# shifts "总计"'s position in the collection. Re-resolve it by name so the
# reference used below tracks the sheet itself rather than its old slot.
$total = $wb.Worksheets.Item("总计")

# Seed the header row + an 8-row data template by copying formatting from
# the structurally-identical "2021-Q4" sheet (gives us the bold/bordered
# header + index-column style used throughout this workbook).
$template.Range("A1:H1").Copy($newSheet.Range("A1"))
$template.Range("A2:H2").Copy($newSheet.Range("A2:H8"))

# Columns B-G hold text (fund code/name/scale/position numbers are kept as
# strings in this workbook, e.g. to preserve leading zeros / fixed decimals).
$newSheet.Range("B2:G8").NumberFormat = "@"

$fundRows = @(
    @(0, "005669", "前海开源公用事业行业股票",                     "258.16", "94.53", "3.33", "8.5967", 10),
    @(1, "008404", "华泰紫金泰盈混合A",                             "4.29",   "79.73", "4.00", "0.1716", 7),
    @(2, "008405", "华泰紫金泰盈混合C",                             "3.71",   "79.73", "4.00", "0.1484", 7),
    @(3, "011694", "华泰紫金信息科技主题6个月定期开放混合A",        "2.60",   "77.49", "4.00", "0.1040", 7),
    @(4, "011695", "华泰紫金信息科技主题6个月定期开放混合C",        "0.83",   "77.49", "4.00", "0.0332", 7),
    @(5, "006923", "前海开源沪港深非周期性行业股票A",               "0.54",   "93.77", "5.22", "0.0282", 7),
    @(6, "006924", "前海开源沪港深非周期性行业股票C",               "0.22",   "93.77", "5.22", "0.0115", 7)
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Add the 2022-Q1 summary row at the top of "总计"
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()

# Drop the stray formatting the row-insert copied down into the data cells.
$total.Range("B2:D2").ClearFormats()

# A2 keeps the bold/bordered index-column style used by the rest of column A.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 9.09

# Renumber the existing index column (previously 0,1,2,3) now that a new
# row sits above them.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
